# simulation-arguments.xlsx — "Switch sim to package"
#
# Renames the cv_obs / sigma_t parameter columns to sd_obs / gp_sigma and
# rescales the gp_scale sweep (column D, and the mirrored column H on the
# gp_scale rows) now that the simulation draws from the new package.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename two of the parameter columns ------------------
$ws.Range("F1").Value = "gp_sigma"
$ws.Range("E1").Value = "sd_obs"

# --- gp_scale column (D) rescaled from the 0.25 baseline to 2 ----------
# (every row except the dedicated gp_scale sweep, rows 15-18)
"2","3","4","5","6","7","8","9","10","11","12","13","14","19","20","21","22","23" | ForEach-Object {
    $ws.Range("D$_").Value = 2
}

# --- gp_scale sweep rows (15-18): values scaled up by 10x ---------------
# column D (gp_scale value) and column H (its mirror) both change together
$ws.Range("D15").Value = 0.5
$ws.Range("H15").Value = 0.5

$ws.Range("D16").Value = 1
$ws.Range("H16").Value = 1

$ws.Range("D17").Value = 4
$ws.Range("H17").Value = 4

$ws.Range("D18").Value = 10
$ws.Range("H18").Value = 10

# --- sd_obs sweep rows (19-23): case label follows the header rename ---
"19","20","21","22","23" | ForEach-Object {
    $ws.Range("G$_").Value = "sd_obs"
}

# --- Selection left where the edit was made -----------------------------
$ws.Range("F2:F22").Select()
